$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the existing X/Y data (A2:B5) with the new series and extend it
# down to row 8 (A2:B8).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1.27

$ws.Range("A3").Value = 0.2
$ws.Range("B3").Value = 2.14

$ws.Range("A4").Value = 0.4
$ws.Range("B4").Value = 3.01

$ws.Range("A5").Value = 0.6
$ws.Range("B5").Value = 3.86

$ws.Range("A6").Value = 0.8
$ws.Range("B6").Value = 3.45

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 2.0499999999999998

$ws.Range("A8").Value = 1.2
$ws.Range("B8").Value = 1.75

# Move the active selection to match the saved view state.
$ws.Range("F23").Select()
